# Auto-generated script applying the value refresh described in the commit diff.
# Updates currentAveragePrice* / Leve* profit columns (H-N) across all 8 server sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 899.4545000000001
$ws.Range("I33").Value = 928.9474
$ws.Range("J33").Value = 712.6667
$ws.Range("K33").Value = 928.9474
$ws.Range("L33").Value = 712.6667
$ws.Range("M33").Value = -699.9474
$ws.Range("N33").Value = -1170.6667
$ws.Range("H64").Value = 3771.4285
$ws.Range("I64").Value = 3600
$ws.Range("J64").Value = 4000
$ws.Range("K64").Value = 3600
$ws.Range("L64").Value = 4000
$ws.Range("M64").Value = -3352
$ws.Range("N64").Value = -4496
$ws.Range("H67").Value = 3771.4285
$ws.Range("I67").Value = 3600
$ws.Range("J67").Value = 4000
$ws.Range("K67").Value = 3600
$ws.Range("L67").Value = 4000
$ws.Range("M67").Value = -2742
$ws.Range("N67").Value = -5716
$ws.Range("H116").Value = 6492.933
$ws.Range("J116").Value = 4599.643
$ws.Range("L116").Value = 4599.643
$ws.Range("N116").Value = -11483.643
$ws.Range("H125").Value = 2768.0667
$ws.Range("I125").Value = 794.8570999999999
$ws.Range("J125").Value = 4494.625
$ws.Range("K125").Value = 7153.7139
$ws.Range("L125").Value = 40451.625
$ws.Range("M125").Value = -4693.7139
$ws.Range("N125").Value = -45371.625
$ws.Range("H132").Value = 12413.63
$ws.Range("I132").Value = 8461.375
$ws.Range("K132").Value = 25384.125
$ws.Range("M132").Value = -22854.125
$ws.Range("H135").Value = 651.9783
$ws.Range("I135").Value = 426.68292
$ws.Range("K135").Value = 3840.14628
$ws.Range("M135").Value = -1305.14628

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 382.83334
$ws.Range("J14").Value = 373.25
$ws.Range("L14").Value = 373.25
$ws.Range("N14").Value = -723.25
$ws.Range("H32").Value = 4139.7935
$ws.Range("I32").Value = 3744.2183
$ws.Range("J32").Value = 6859.375
$ws.Range("K32").Value = 3744.2183
$ws.Range("L32").Value = 6859.375
$ws.Range("M32").Value = -3457.2183
$ws.Range("N32").Value = -7433.375
$ws.Range("H63").Value = 3436.25
$ws.Range("I63").Value = 3099.4443
$ws.Range("J63").Value = 3711.818
$ws.Range("K63").Value = 3099.4443
$ws.Range("L63").Value = 3711.818
$ws.Range("M63").Value = -2413.4443
$ws.Range("N63").Value = -5083.818
$ws.Range("H66").Value = 3436.25
$ws.Range("I66").Value = 3099.4443
$ws.Range("J66").Value = 3711.818
$ws.Range("K66").Value = 15497.2215
$ws.Range("L66").Value = 18559.09
$ws.Range("M66").Value = -12065.2215
$ws.Range("N66").Value = -25423.09
$ws.Range("H122").Value = 3474.111
$ws.Range("I122").Value = 3452.3076
$ws.Range("J122").Value = 3494.3572
$ws.Range("K122").Value = 10356.9228
$ws.Range("L122").Value = 10483.0716
$ws.Range("M122").Value = -7906.9228
$ws.Range("N122").Value = -15383.0716

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 70447.5
$ws.Range("J40").Value = 70447.5
$ws.Range("L40").Value = 70447.5
$ws.Range("N40").Value = -70977.5
$ws.Range("H82").Value = 20901.375
$ws.Range("I82").Value = 16744.572
$ws.Range("K82").Value = 16744.572
$ws.Range("M82").Value = -16361.572
$ws.Range("H85").Value = 20901.375
$ws.Range("I85").Value = 16744.572
$ws.Range("K85").Value = 16744.572
$ws.Range("M85").Value = -15418.572
$ws.Range("H96").Value = 8000
$ws.Range("I96").Value = 8000
$ws.Range("K96").Value = 8000
$ws.Range("M96").Value = -5254
$ws.Range("H107").Value = 1864.8334
$ws.Range("I107").Value = 1961.5
$ws.Range("K107").Value = 1961.5
$ws.Range("M107").Value = -41.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1423.8823
$ws.Range("J16").Value = 1631.5834
$ws.Range("L16").Value = 1631.5834
$ws.Range("N16").Value = -2205.5834
$ws.Range("H31").Value = 2153.2942
$ws.Range("J31").Value = 2936
$ws.Range("L31").Value = 2936
$ws.Range("N31").Value = -3526
$ws.Range("H34").Value = 2153.2942
$ws.Range("J34").Value = 2936
$ws.Range("L34").Value = 2936
$ws.Range("N34").Value = -3340
$ws.Range("H99").Value = 7023.6265
$ws.Range("I99").Value = 8558.5625
$ws.Range("K99").Value = 8558.5625
$ws.Range("M99").Value = -7060.5625
$ws.Range("H113").Value = 1423.8823
$ws.Range("J113").Value = 1631.5834
$ws.Range("L113").Value = 1631.5834
$ws.Range("N113").Value = -5971.5834
$ws.Range("H122").Value = 4344.243
$ws.Range("I122").Value = 4275.952
$ws.Range("J122").Value = 4433.875
$ws.Range("K122").Value = 12827.856
$ws.Range("L122").Value = 13301.625
$ws.Range("M122").Value = -10377.856
$ws.Range("N122").Value = -18201.625
$ws.Range("H126").Value = 7023.6265
$ws.Range("I126").Value = 8558.5625
$ws.Range("K126").Value = 25675.6875
$ws.Range("M126").Value = -23205.6875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 76760.375
$ws.Range("J37").Value = 76760.375
$ws.Range("L37").Value = 230281.125
$ws.Range("N37").Value = -230505.125
$ws.Range("H82").Value = 17221.25
$ws.Range("I82").Value = 14443.5
$ws.Range("K82").Value = 43330.5
$ws.Range("M82").Value = -42924.5
$ws.Range("H85").Value = 17221.25
$ws.Range("I85").Value = 14443.5
$ws.Range("K85").Value = 43330.5
$ws.Range("M85").Value = -41926.5
$ws.Range("H116").Value = 3924021.2
$ws.Range("I116").Value = 3924021.2
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 11772063.6
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -11768621.6
$ws.Range("N116").ClearContents()
$ws.Range("H128").Value = 180041.67
$ws.Range("I128").Value = 180041.67
$ws.Range("K128").Value = 540125.01
$ws.Range("M128").Value = -535145.01
$ws.Range("H137").Value = 11580.632
$ws.Range("I137").Value = 6316.1665
$ws.Range("J137").Value = 14010.385
$ws.Range("K137").Value = 18948.4995
$ws.Range("L137").Value = 42031.155
$ws.Range("M137").Value = -13848.4995
$ws.Range("N137").Value = -52231.155

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 64962.5
$ws.Range("I41").Value = 64962.5
$ws.Range("K41").Value = 64962.5
$ws.Range("M41").Value = -64607.5
$ws.Range("H70").Value = 6312.933
$ws.Range("J70").Value = 5857.9165
$ws.Range("L70").Value = 5857.9165
$ws.Range("N70").Value = -6397.9165
$ws.Range("H73").Value = 6312.933
$ws.Range("J73").Value = 5857.9165
$ws.Range("L73").Value = 5857.9165
$ws.Range("N73").Value = -7729.9165
$ws.Range("H113").Value = 2729.4
$ws.Range("I113").Value = 1361.8334
$ws.Range("K113").Value = 1361.8334
$ws.Range("M113").Value = 808.1666
$ws.Range("H126").Value = 7333.4287
$ws.Range("I126").Value = 6152.0713
$ws.Range("J126").Value = 9696.143
$ws.Range("K126").Value = 18456.2139
$ws.Range("L126").Value = 29088.429
$ws.Range("M126").Value = -15986.2139
$ws.Range("N126").Value = -34028.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2588.5
$ws.Range("I7").Value = 2000
$ws.Range("K7").Value = 2000
$ws.Range("M7").Value = -1888
$ws.Range("H22").Value = 820.05884
$ws.Range("I22").Value = 823.5454999999999
$ws.Range("K22").Value = 823.5454999999999
$ws.Range("M22").Value = -528.5454999999999
$ws.Range("H27").Value = 820.05884
$ws.Range("I27").Value = 823.5454999999999
$ws.Range("K27").Value = 823.5454999999999
$ws.Range("M27").Value = -716.5454999999999
$ws.Range("H40").Value = 4583.1113
$ws.Range("I40").Value = 4469.12
$ws.Range("K40").Value = 4469.12
$ws.Range("M40").Value = -4333.12
$ws.Range("H61").Value = 17576.666
$ws.Range("I61").Value = 17576.666
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 17576.666
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -17374.666
$ws.Range("N61").ClearContents()
$ws.Range("H105").Value = 66271.5
$ws.Range("J105").Value = 66271.5
$ws.Range("L105").Value = 66271.5
$ws.Range("N105").Value = -73259.5
$ws.Range("H113").Value = 17576.666
$ws.Range("I113").Value = 17576.666
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 17576.666
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -15406.666
$ws.Range("N113").ClearContents()
$ws.Range("H126").Value = 2588.5
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 8250
$ws.Range("I33").Value = 8250
$ws.Range("K33").Value = 8250
$ws.Range("M33").Value = -8000
$ws.Range("H36").Value = 8250
$ws.Range("I36").Value = 8250
$ws.Range("K36").Value = 8250
$ws.Range("M36").Value = -8000
$ws.Range("H37").Value = 18999.834
$ws.Range("I37").Value = 19249.75
$ws.Range("K37").Value = 19249.75
$ws.Range("M37").Value = -19046.75
$ws.Range("H41").Value = 12198.5
$ws.Range("I41").Value = 11364.667
$ws.Range("K41").Value = 11364.667
$ws.Range("M41").Value = -10974.667
$ws.Range("H105").Value = 53500
$ws.Range("J105").Value = 53500
$ws.Range("L105").Value = 53500
$ws.Range("N105").Value = -60488
$ws.Range("H122").Value = 19394.857
$ws.Range("I122").Value = 7153.2
$ws.Range("K122").Value = 21459.6
$ws.Range("M122").Value = -19009.6
$ws.Range("H126").Value = 4888.8
$ws.Range("I126").Value = 4504.2354
$ws.Range("K126").Value = 13512.7062
$ws.Range("M126").Value = -11042.7062

